$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update temporal_freq values (B2) from "[1, 5]" to "[0.5, 1, 2]"
$ws.Range("B2").Value = "[0.5, 1, 2]"

# Update repetitions (E2) from 5 to 1
$ws.Range("E2").Value = 1

# Move the active selection from F2 to C9
$ws.Range("C9").Select()
